$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 241; existing rows 241-302 shift down to 242-303.
$ws.Rows.Item(241).Insert()

# Populate the newly inserted row 241 with a new price record
# (same dimension values as the former row 241, but new date/volume/price figures).
$ws.Cells.Item(241, 1).Value = 10
$ws.Cells.Item(241, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(241, 3).Value = "La Araucanía"
$ws.Cells.Item(241, 4).Value = 45204
$ws.Cells.Item(241, 5).Value = 9
$ws.Cells.Item(241, 6).Value = "Fruta"
$ws.Cells.Item(241, 7).Value = 100104
$ws.Cells.Item(241, 8).Value = "Frutos de pepita"
$ws.Cells.Item(241, 9).Value = 100104001
$ws.Cells.Item(241, 10).Value = "Granada"
$ws.Cells.Item(241, 11).Value = "Wonderfull"
$ws.Cells.Item(241, 12).Value = "Primera"
$ws.Cells.Item(241, 13).Value = 400
$ws.Cells.Item(241, 14).Value = 17000
$ws.Cells.Item(241, 15).Value = 17000
$ws.Cells.Item(241, 16).Value = 17000
$ws.Cells.Item(241, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(241, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(241, 19).Value = 1700
$ws.Cells.Item(241, 20).Value = 10
